# "Se procesan de nuevo los datos con las nuevas dimensiones curadas"
# Re-process the metadata sheet with the newly curated dimensions:
# - municipio-nombre (col C) becomes a refArea dimension instead of a measure.
# - condicion-socioeconomica (col E) and sexo (col G) become measures instead
#   of skos:Concept dimensions (no more external mapping workbook needed).
# - aragon (col F) keeps its dim role but now points at a Comunidad URI
#   instead of a generic skos:Concept.
# - The old row of "mapping-*.xlsx" helper references (row 5) is removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: predicate ---
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:condicion-socioeconomica"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:sexo"

# --- Row 3: dim / medida role ---
$ws.Range("C3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("G3").Value = "medida"

# --- Row 4: datatype / URI template ---
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "xsd:int"

# --- Row 5: the mapping-*.xlsx helper row is no longer needed ---
$ws.Range("A5:K5").EntireRow.Delete()
